# evidence updates & sheet refactors
# Remove the "delta_duc" column (column D) from the sheet/table, and
# re-apply the author's manual column-width tweaks + selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet's single table (Table1) includes the column being removed,
# so drop the table's structure first, delete the worksheet column, and
# rebuild the table over the new (narrower) range so its column list and
# header names stay in sync with the sheet.
$lo = $ws.ListObjects.Item(1)
$lo.Unlist()

$ws.Range("D1").EntireColumn.Delete()

$lo2 = $ws.ListObjects.Add(1, $ws.Range("A1:H2"), $null, 1)
$lo2.Name = "Table1"
$lo2.TableStyle = "TableStyleMedium15"

# Manual column-width adjustments made after the deletion (target raw
# widths ~90.1640625 / 19.5 pt; ColumnWidth is in characters, so nudge to
# the closest value that rounds to the desired stored width).
$ws.Columns("D:D").ColumnWidth = 89.33
$ws.Columns("F:F").ColumnWidth = 18.67

# Final selection left on the (now widened) delta_pattern column.
$ws.Columns("D:D").Select() | Out-Null
